# March 16 covid data update:
# - "Covid5" header (and its data) moves from column P to column O
#   (column O was previously empty, so this is effectively adding the
#   "Covid5" data column right after "Covid4").
# - Selection moves to H11.
# - Page orientation is set to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Covid5" header from P1 into O1, then clear the old P1 cell.
$ws.Range("O1").Value = $ws.Range("P1").Text
$ws.Range("P1").ClearContents()

# New "Covid5" values for each county row (row 2 .. row 59).
$covid5 = @(
    19, 0, 5, 0, 0, 0, 6, 0, 1, 1, 0, 2, 37, 12, 0, 0, 0, 2, 1, 0, 2, 1, 0, 0,
    2, 21, 1, 629, 2, 4, 18, 2, 113, 0, 2, 1, 10, 1, 229, 1, 9, 7, 0, 0, 0, 3,
    475, 14, 1, 6, 13, 0, 0, 1, 1253, 3, 0, 3537
)

for ($i = 0; $i -lt $covid5.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 15).Value = $covid5[$i]
}

# Page setup: portrait orientation.
$ws.PageSetup.Orientation = 1

# Update the active selection to H11.
$ws.Range("H11").Select()
